$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-formatted cells need an explicit Text number format before the
# assignment, otherwise Excel's "smart" entry parses a trailing "%" into a
# percentage number (e.g. 63% -> 0.63) instead of keeping the literal string.
$percentCells = @("H13","H14","H15","H17","H21","H22","H26","H29","H32","H34","H35","H41","H45")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-26 18:48:20"
$ws.Range("E3").Value = "2026-02-26 18:48:23"
$ws.Range("E4").Value = "2026-02-26 18:48:25"
$ws.Range("E5").Value = "2026-02-26 18:48:28"
$ws.Range("E6").Value = "2026-02-26 18:48:30"
$ws.Range("E7").Value = "2026-02-26 18:48:33"
$ws.Range("E8").Value = "2026-02-26 18:48:35"
$ws.Range("O8").Value = "11.4 °C"
$ws.Range("E9").Value = "2026-02-26 18:48:37"
$ws.Range("E10").Value = "2026-02-26 18:48:38"
$ws.Range("E11").Value = "2026-02-26 18:48:39"
$ws.Range("O11").Value = "8.9 °C"
$ws.Range("E12").Value = "2026-02-26 18:48:40"
$ws.Range("O12").Value = "11.6 °C"
$ws.Range("E13").Value = "2026-02-26 18:48:42"
$ws.Range("H13").Value = "63%"
$ws.Range("J13").Value = "1028.2 hPa"
$ws.Range("O13").Value = "7.2 °C"
$ws.Range("E14").Value = "2026-02-26 18:48:43"
$ws.Range("H14").Value = "87%"
$ws.Range("O14").Value = "11.9 °C"
$ws.Range("E15").Value = "2026-02-26 18:48:44"
$ws.Range("H15").Value = "83%"
$ws.Range("O15").Value = "12.1 °C"
$ws.Range("E16").Value = "2026-02-26 18:48:45"
$ws.Range("E17").Value = "2026-02-26 18:48:46"
$ws.Range("H17").Value = "38%"
$ws.Range("E18").Value = "2026-02-26 18:48:47"
$ws.Range("J18").Value = "1027.2 hPa"
$ws.Range("E19").Value = "2026-02-26 18:48:48"
$ws.Range("E20").Value = "2026-02-26 18:48:51"
$ws.Range("O20").Value = "2.9 °C"
$ws.Range("E21").Value = "2026-02-26 18:48:53"
$ws.Range("H21").Value = "62%"
$ws.Range("J21").Value = "1027.0 hPa"
$ws.Range("O21").Value = "10.0 °C"
$ws.Range("E22").Value = "2026-02-26 18:48:56"
$ws.Range("H22").Value = "49%"
$ws.Range("E23").Value = "2026-02-26 18:48:58"
$ws.Range("E24").Value = "2026-02-26 18:49:00"
$ws.Range("E25").Value = "2026-02-26 18:49:03"
$ws.Range("O25").Value = "5.5 °C"
$ws.Range("E26").Value = "2026-02-26 18:49:05"
$ws.Range("H26").Value = "39%"
$ws.Range("O26").Value = "11.2 °C"
$ws.Range("E27").Value = "2026-02-26 18:49:08"
$ws.Range("E28").Value = "2026-02-26 18:49:10"
$ws.Range("O28").Value = "11.2 °C"
$ws.Range("E29").Value = "2026-02-26 18:49:13"
$ws.Range("H29").Value = "85%"
$ws.Range("E30").Value = "2026-02-26 18:49:15"
$ws.Range("E31").Value = "2026-02-26 18:49:18"
$ws.Range("O31").Value = "12.0 °C"
$ws.Range("E32").Value = "2026-02-26 18:49:20"
$ws.Range("H32").Value = "62%"
$ws.Range("O32").Value = "8.5 °C"
$ws.Range("E33").Value = "2026-02-26 18:49:23"
$ws.Range("O33").Value = "8.7 °C"
$ws.Range("E34").Value = "2026-02-26 18:49:25"
$ws.Range("H34").Value = "45%"
$ws.Range("O34").Value = "5.1 °C"
$ws.Range("E35").Value = "2026-02-26 18:49:28"
$ws.Range("H35").Value = "41%"
$ws.Range("E36").Value = "2026-02-26 18:49:30"
$ws.Range("E37").Value = "2026-02-26 18:49:33"
$ws.Range("E38").Value = "2026-02-26 18:49:35"
$ws.Range("E39").Value = "2026-02-26 18:49:38"
$ws.Range("O39").Value = "3.0 °C"
$ws.Range("E40").Value = "2026-02-26 18:49:40"
$ws.Range("E41").Value = "2026-02-26 18:49:43"
$ws.Range("H41").Value = "82%"
$ws.Range("E42").Value = "2026-02-26 18:49:45"
$ws.Range("E43").Value = "2026-02-26 18:49:48"
$ws.Range("O43").Value = "9.5 °C"
$ws.Range("E44").Value = "2026-02-26 18:49:50"
$ws.Range("E45").Value = "2026-02-26 18:49:52"
$ws.Range("G45").Value = "2 cm"
$ws.Range("H45").Value = "46%"
$ws.Range("E46").Value = "2026-02-26 18:49:55"
$ws.Range("O46").Value = "11.3 °C"
